$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 334
$ws.Range("H2").Value = 296
$ws.Range("I2").Value = 1284
$ws.Range("J2").Value = 1284
$ws.Range("K2").Value = 1284
$ws.Range("L2").Value = 1284
$ws.Range("M2").Value = 1284
$ws.Range("N2").Value = 1284
$ws.Range("O2").Value = 1284
$ws.Range("P2").Value = 1284
$ws.Range("Q2").Value = 1284
$ws.Range("R2").Value = 1284
$ws.Range("S2").Value = 1284
$ws.Range("T2").Value = 1284
$ws.Range("U2").Value = 1284
$ws.Range("V2").Value = 1284
$ws.Range("W2").Value = 1284
$ws.Range("X2").Value = 1284
$ws.Range("Y2").Value = 1284
$ws.Range("Z2").Value = 1284
$ws.Range("AA2").Value = 1284
$ws.Range("AB2").Value = 1284
$ws.Range("AC2").Value = 1284
$ws.Range("AD2").Value = 1284
$ws.Range("AE2").Value = 1284
$ws.Range("AF2").Value = 1284
$ws.Range("AG2").Value = 1284
$ws.Range("AH2").Value = 1284
$ws.Range("AI2").Value = 1284
$ws.Range("AJ2").Value = 1284
$ws.Range("AK2").Value = 1284
$ws.Range("AL2").Value = 1284
$ws.Range("AM2").Value = 1284
$ws.Range("AN2").Value = 1284
$ws.Range("AO2").Value = 1284
$ws.Range("AP2").Value = 1284
$ws.Range("AQ2").Value = 1284
$ws.Range("AR2").Value = 1284
$ws.Range("AS2").Value = 1284
$ws.Range("AT2").Value = 1284
$ws.Range("AU2").Value = 1284
$ws.Range("AV2").Value = 1284
$ws.Range("AW2").Value = 1284
$ws.Range("AX2").Value = 1284
$ws.Range("AY2").Value = 1284
$ws.Range("AZ2").Value = 1284
$ws.Range("BA2").Value = 1284
$ws.Range("BB2").Value = 1284
$ws.Range("BC2").Value = 1284
$ws.Range("BD2").Value = 1284
$ws.Range("BE2").Value = 1284
$ws.Range("BF2").Value = 1284
$ws.Range("BG2").Value = 1284
$ws.Range("BH2").Value = 1284
$ws.Range("BI2").Value = 1284
$ws.Range("C3").Value = 625.9098861047836
$ws.Range("D3").Value = 677.0424634146341
$ws.Range("E3").Value = 1076.764257668712
$ws.Range("F3").Value = 1548.690778443114
$ws.Range("G3").Value = 94.64796992481202
$ws.Range("H3").Value = 9.976722972972972
$ws.Range("I3").Value = 7.675436137071651
$ws.Range("J3").Value = 1.737461059190031
$ws.Range("K3").Value = 2.501658878504673
$ws.Range("L3").Value = 0.8867757009345796
$ws.Range("M3").Value = 1.946246105919003
$ws.Range("N3").Value = 0.7805062305295951
$ws.Range("O3").Value = 10.05654205607476
$ws.Range("P3").Value = 0.719158878504673
$ws.Range("Q3").Value = 1.701191588785047
$ws.Range("R3").Value = 0.8104361370716511
$ws.Range("S3").Value = 3.335303738317758
$ws.Range("T3").Value = 0.4047507788161994
$ws.Range("U3").Value = 0.03353582554517134
$ws.Range("V3").Value = 2.043146417445483
$ws.Range("W3").Value = 14.49653426791278
$ws.Range("X3").Value = 16.01971183800623
$ws.Range("Y3").Value = 4.019369158878505
$ws.Range("Z3").Value = 6.014563862928348
$ws.Range("AA3").Value = 5.813971962616822
$ws.Range("AB3").Value = 1.571339563862928
$ws.Range("AC3").Value = 0.7404205607476636
$ws.Range("AD3").Value = 3.504361370716511
$ws.Range("AE3").Value = 0.04205607476635514
$ws.Range("AF3").Value = 1.111020249221184
$ws.Range("AG3").Value = 0.1230996884735202
$ws.Range("AH3").Value = 0.007009345794392523
$ws.Range("AI3").Value = 3.732554517133956
$ws.Range("AJ3").Value = 0.5957943925233645
$ws.Range("AK3").Value = 0.3642523364485981
$ws.Range("AL3").Value = 0.1582788161993769
$ws.Range("AM3").Value = 0.1580996884735202
$ws.Range("AN3").Value = 0.9637850467289719
$ws.Range("AO3").Value = 0.1627725856697819
$ws.Range("AP3").Value = 2.03714953271028
$ws.Range("AQ3").Value = 0.1362928348909657
$ws.Range("AR3").Value = 0.839151090342679
$ws.Range("AS3").Value = 0.1515809968847352
$ws.Range("AT3").Value = 0.1929828660436137
$ws.Range("AU3").Value = 0.270404984423676
$ws.Range("AV3").Value = 0.995880062305296
$ws.Range("AW3").Value = 0.3393380062305296
$ws.Range("AX3").Value = 0.09959501557632398
$ws.Range("AY3").Value = 0.07228193146417446
$ws.Range("AZ3").Value = 0.07398753894080996
$ws.Range("BA3").Value = 0.1082554517133956
$ws.Range("BB3").Value = 0.1094236760124611
$ws.Range("BC3").Value = 0.001557632398753894
$ws.Range("BD3").Value = 0.2414330218068536
$ws.Range("BE3").Value = 0.08566978193146417
$ws.Range("BH3").Value = 0.01129283489096573
$ws.Range("BI3").Value = 0.8403426791277259
$ws.Range("C4").Value = 171.6169439449899
$ws.Range("D4").Value = 175.981797166909
$ws.Range("E4").Value = 265.5547198582155
$ws.Range("F4").Value = 495.0505359133559
$ws.Range("G4").Value = 52.5145269745037
$ws.Range("H4").Value = 12.51796397606443
$ws.Range("I4").Value = 14.36218134157754
$ws.Range("J4").Value = 5.259931890642274
$ws.Range("K4").Value = 6.720174751292573
$ws.Range("L4").Value = 3.67738093152069
$ws.Range("M4").Value = 5.670064484618109
$ws.Range("N4").Value = 2.825258825071417
$ws.Range("O4").Value = 21.47763890813834
$ws.Range("P4").Value = 2.96830042245807
$ws.Range("Q4").Value = 5.216002273481808
$ws.Range("R4").Value = 3.771255123065122
$ws.Range("S4").Value = 8.201720138837237
$ws.Range("T4").Value = 5.064132815135191
$ws.Range("U4").Value = 0.6633434922714724
$ws.Range("V4").Value = 11.11792488952686
$ws.Range("W4").Value = 16.58521717973746
$ws.Range("X4").Value = 23.05497046131782
$ws.Range("Y4").Value = 9.126271789740716
$ws.Range("Z4").Value = 7.901806706384862
$ws.Range("AA4").Value = 16.49482156729849
$ws.Range("AB4").Value = 6.379426126646438
$ws.Range("AC4").Value = 5.012613254485819
$ws.Range("AD4").Value = 13.8277850285873
$ws.Range("AE4").Value = 0.4254464181144004
$ws.Range("AF4").Value = 5.22427120624288
$ws.Range("AG4").Value = 0.829359199733317
$ws.Range("AH4").Value = 0.2031262213472014
$ws.Range("AI4").Value = 13.3145025280377
$ws.Range("AJ4").Value = 3.411127233352597
$ws.Range("AK4").Value = 4.230644723870419
$ws.Range("AL4").Value = 2.345752516447703
$ws.Range("AM4").Value = 2.35489337517622
$ws.Range("AN4").Value = 7.634456800255992
$ws.Range("AO4").Value = 2.549124607830517
$ws.Range("AP4").Value = 11.05991469889155
$ws.Range("AQ4").Value = 2.131153463473378
$ws.Range("AR4").Value = 4.019023001713187
$ws.Range("AS4").Value = 1.501242636902157
$ws.Range("AT4").Value = 2.858069332321699
$ws.Range("AU4").Value = 3.948712053549328
$ws.Range("AV4").Value = 4.470390898140147
$ws.Range("AW4").Value = 3.082851234817126
$ws.Range("AX4").Value = 1.880549783918487
$ws.Range("AY4").Value = 0.9329377722750313
$ws.Range("AZ4").Value = 1.89719635288587
$ws.Range("BA4").Value = 2.26050030884912
$ws.Range("BB4").Value = 2.054380710595965
$ws.Range("BC4").Value = 0.05581455721859475
$ws.Range("BD4").Value = 3.489786429629616
$ws.Range("BE4").Value = 1.771114977914812
$ws.Range("BH4").Value = 0.1734215168908341
$ws.Range("BI4").Value = 0.3664306667457159
$ws.Range("F6").Value = 1312.75
$ws.Range("H6").Value = 1.6225
$ws.Range("H7").Value = 5.52
$ws.Range("F8").Value = 1797.875
$ws.Range("H8").Value = 11
$ws.Range("O8").Value = 3
$ws.Range("W8").Value = 25
$ws.Range("X8").Value = 35
$ws.Range("J9").Value = 27.27
$ws.Range("P9").Value = 23.08
$ws.Range("S9").Value = 64.8
$ws.Range("AV9").Value = 42.86
$ws.Range("BH9").Value = 4

Write-Host "Applied changes"